$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row 25 with new task first
$ws.Range("A25").Value = "sistematizar el RG 14"
$ws.Range("B25").Value = "no comenzado"

# Update row 23: append text to task and change status to "en proceso"
$ws.Range("A23").Value = "relacionar ordenes a la inversa, revisar porque hacia abajo trae 2 veces el mismo item"
$ws.Range("B23").Value = "en proceso"

$ws.Range("C25").Select()
